$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "FINAL PROJECT" heading (paragraph 2): grow to sz=32 and append " - DAG".
#    First insert a NEW blank paragraph *after* it (while it still carries
#    the old sz=30 centered/bold formatting) so the blank paragraph keeps
#    that original look, matching the target's spacer paragraph.
# ---------------------------------------------------------------------------
$pFinalProject = $d.Paragraphs(2)
$pFinalProject.Range.InsertParagraphAfter()

# Re-fetch paragraph 2 (still "FINAL PROJECT") and bump its size to 32
# (half-points), i.e. 16pt, across the whole paragraph (incl. the paragraph
# mark) so the pPr/rPr mark formatting updates too.
$pFinalProject = $d.Paragraphs(2)
$pFinalProject.Range.Font.Size = 16

# Append " - DAG" right after the existing text (exclude the trailing
# paragraph mark from the insertion point).
$rFP = $pFinalProject.Range
$rFP.MoveEnd(1, -1)
$rFP.Collapse(0)
$rDag = $d.Range($rFP.End, $rFP.End)
$rDag.Text = " - DAG"
$rDag.Font.Bold = $true
$rDag.Font.Size = 16

# ---------------------------------------------------------------------------
# 2. "Introduction:" paragraph (now paragraph 4): bold, sz=32, underline.
# ---------------------------------------------------------------------------
$pIntro = $d.Paragraphs(4)
$pIntro.Range.Font.Size = 16
$pIntro.Range.Font.Bold = $true
$pIntro.Range.Font.Underline = 1

# ---------------------------------------------------------------------------
# 3. Body paragraph (now paragraph 5): justify + extend its text.
# ---------------------------------------------------------------------------
$pBody = $d.Paragraphs(5)
$pBody.Range.ParagraphFormat.Alignment = 3

$rBody = $pBody.Range
$rBody.MoveEnd(1, -1)
$rBody.Collapse(0)
$rBody.InsertAfter(". In DAG each node is represented with some ID and its VALUE. In that, we count the occurrence of a given number at node value in DAG. Multiple nodes can have same or different value. All the nodes are visited in Breadth-first manner.")

# ---------------------------------------------------------------------------
# 4. Detach the _GoBack bookmark now (it currently sits glued to the end of
#    paragraph 5's original text) so later paragraph insertions don't drag
#    it along; it gets re-created in its final home further below.
# ---------------------------------------------------------------------------
$bmOld = $d.Bookmarks("_GoBack")
$bmOld.Delete()

# ---------------------------------------------------------------------------
# 5. Create all five new trailing paragraphs (empty shells first), then fill
#    in text/formatting explicitly on each so nothing is left over-inherited
#    from its predecessor (e.g. the body's jc=both must NOT leak onto the
#    "Files Structure" paragraphs).
# ---------------------------------------------------------------------------
$d.Paragraphs(5).Range.InsertParagraphAfter()   # 6: "I performed this..."
$d.Paragraphs(6).Range.InsertParagraphAfter()   # 7: "While coding..."
$d.Paragraphs(7).Range.InsertParagraphAfter()   # 8: bookmark-only paragraph
$d.Paragraphs(8).Range.InsertParagraphAfter()   # 9: "Files Structure:"
$d.Paragraphs(9).Range.InsertParagraphAfter()   # 10: "The code zip file..."

# --- paragraph 6: "I performed this using programming in C++. ..." ---
$p6 = $d.Paragraphs(6)
$p6.Range.ParagraphFormat.Alignment = 3
$r6 = $p6.Range
$r6.MoveEnd(1, -1)
$r6.InsertAfter("`tI performed this using programming in C++. Starting with sequential approach and parallelizing it later. The same solution is also done using FastFlow library. ")

# --- paragraph 7: "While coding the parallel form, ..." ---
$p7 = $d.Paragraphs(7)
$p7.Range.ParagraphFormat.Alignment = 3
$r7 = $p7.Range
$r7.MoveEnd(1, -1)
$r7.InsertAfter("`tWhile coding the parallel form, I considered the FARM design pattern with some sort of JOB Stealing mechanism that we studied in the course. A queue is shared with all the threads workers and data is popped and pushed into it for processing and the visited queue take the record of the nodes that have been visited already. This is the overview of the internal working of the system and will be detailed in later sections.")

# --- paragraph 8: bookmark-only paragraph ---
# NB: Bookmarks.Add() on a truly collapsed (Start==End) Range that lies in
# newly-inserted territory is unreliable in this host, so we temporarily
# type a unique placeholder, bookmark the (non-collapsed) range around it,
# then delete the placeholder text -- the bookmark collapses along with it
# and ends up exactly where we want, with no leftover run.
$p8 = $d.Paragraphs(8)
$p8.Range.ParagraphFormat.Alignment = 3
$r8 = $p8.Range
$r8.MoveEnd(1, -1)
$r8.InsertAfter("ZZ_BOOKMARK_PLACEHOLDER_ZZ")

$p8 = $d.Paragraphs(8)
$rFind = $p8.Range
$rFind.Find.Execute("ZZ_BOOKMARK_PLACEHOLDER_ZZ")
$d.Bookmarks.Add("_GoBack", $rFind)

$p8 = $d.Paragraphs(8)
$rFind2 = $p8.Range
$rFind2.Find.Execute("ZZ_BOOKMARK_PLACEHOLDER_ZZ")
$rFind2.Text = ""

# --- paragraph 9: "Files Structure:" heading ---
$p9 = $d.Paragraphs(9)
$p9.Range.ParagraphFormat.Alignment = 0
$p9.Range.Font.Size = 16
$p9.Range.Font.Bold = $true
$p9.Range.Font.Underline = 1
$r9 = $p9.Range
$r9.MoveEnd(1, -1)
$r9.InsertAfter("Files Structure:")

# --- paragraph 10: "The code zip file or github repo ..." ---
$p10 = $d.Paragraphs(10)
$p10.Range.ParagraphFormat.Alignment = 0
$r10 = $p10.Range
$r10.MoveEnd(1, -1)
$r10.InsertAfter("`tThe code zip file or github repo contains the code of this project. ")

foreach ($p in $d.Paragraphs) {
    Write-Output ("[" + $p.Range.Text + "]")
}
